$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 714.5833
$ws.Range("J17").Value = 714.5833
$ws.Range("L17").Value = 2143.7499
$ws.Range("N17").Value = -2479.7499

$ws.Range("H33").Value = 176.32
$ws.Range("I33").Value = 176.32
$ws.Range("K33").Value = 176.32
$ws.Range("M33").Value = 52.68000000000001

$ws.Range("H40").Value = 800
$ws.Range("I40").Value = 800
$ws.Range("K40").Value = 800
$ws.Range("M40").Value = -625

$ws.Range("H92").Value = 1641.5
$ws.Range("I92").Value = 1661.6154
$ws.Range("J92").Value = 1380
$ws.Range("K92").Value = 1661.6154
$ws.Range("L92").Value = 1380
$ws.Range("M92").Value = -413.6153999999999
$ws.Range("N92").Value = -3876

$ws.Range("H98").Value = 899
$ws.Range("I98").Value = 959
$ws.Range("J98").Value = 599
$ws.Range("K98").Value = 959
$ws.Range("L98").Value = 599
$ws.Range("M98").Value = 539
$ws.Range("N98").Value = -3595

$ws.Range("H116").Value = 4075.9092
$ws.Range("I116").Value = 4095.2222
$ws.Range("J116").Value = 3989
$ws.Range("K116").Value = 4095.2222
$ws.Range("L116").Value = 3989
$ws.Range("M116").Value = -653.2222000000002
$ws.Range("N116").Value = -10873

$ws.Range("H122").Value = 899
$ws.Range("I122").Value = 959
$ws.Range("J122").Value = 599
$ws.Range("K122").Value = 2877
$ws.Range("L122").Value = 1797
$ws.Range("M122").Value = -427
$ws.Range("N122").Value = -6697

$ws.Range("H131").Value = 9641
$ws.Range("I131").Value = 9775
$ws.Range("K131").Value = 29325
$ws.Range("M131").Value = -24285

$ws.Range("H135").Value = 3162.3333
$ws.Range("I135").Value = 4037.6667
$ws.Range("K135").Value = 36339.0003
$ws.Range("M135").Value = -33804.0003

$ws.Range("H137").Value = 2944.818
$ws.Range("I137").Value = 1770.5714
$ws.Range("K137").Value = 5311.7142
$ws.Range("M137").Value = -2761.7142

$ws.Range("H138").Value = 6983.2812
$ws.Range("J138").Value = 7684.963
$ws.Range("L138").Value = 23054.889
$ws.Range("N138").Value = -33334.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2951.75
$ws.Range("I45").Value = 2546.1428
$ws.Range("J45").Value = 3519.6
$ws.Range("K45").Value = 2546.1428
$ws.Range("L45").Value = 3519.6
$ws.Range("M45").Value = -2169.1428
$ws.Range("N45").Value = -4273.6

$ws.Range("H61").Value = 3471
$ws.Range("I61").Value = 3285.75
$ws.Range("K61").Value = 3285.75
$ws.Range("M61").Value = -3073.75

$ws.Range("H74").Value = 1508.8889
$ws.Range("I74").Value = 797.1429
$ws.Range("K74").Value = 797.1429
$ws.Range("M74").Value = 76.85709999999995

$ws.Range("H77").Value = 1508.8889
$ws.Range("I77").Value = 797.1429
$ws.Range("K77").Value = 3985.7145
$ws.Range("M77").Value = 382.2855

$ws.Range("H97").Value = 798
$ws.Range("I97").Value = 798
$ws.Range("K97").Value = 798
$ws.Range("M97").Value = -302

$ws.Range("H122").Value = 3081.8333
$ws.Range("I122").Value = 3495.4
$ws.Range("K122").Value = 10486.2
$ws.Range("M122").Value = -8036.200000000001

$ws.Range("H132").Value = 3522.5625
$ws.Range("I132").Value = 2766.7144
$ws.Range("J132").Value = 4110.4443
$ws.Range("K132").Value = 8300.143199999999
$ws.Range("L132").Value = 12331.3329
$ws.Range("M132").Value = -5770.143199999999
$ws.Range("N132").Value = -17391.3329

$ws.Range("H136").Value = 3471
$ws.Range("I136").Value = 3285.75
$ws.Range("K136").Value = 9857.25
$ws.Range("M136").Value = -7307.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 66666
$ws.Range("J57").Value = 66666
$ws.Range("L57").Value = 66666
$ws.Range("N57").Value = -68106

$ws.Range("H134").Value = 3669.842
$ws.Range("I134").Value = 3651.5
$ws.Range("K134").Value = 10954.5
$ws.Range("M134").Value = -8419.5

$ws.Range("H136").Value = 66666
$ws.Range("J136").Value = 66666
$ws.Range("L136").Value = 66666
$ws.Range("N136").Value = -76866

$ws.Range("H137").Value = 66666
$ws.Range("J137").Value = 66666
$ws.Range("L137").Value = 66666
$ws.Range("N137").Value = -76866

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4834.846
$ws.Range("I31").Value = 2872.7058
$ws.Range("K31").Value = 2872.7058
$ws.Range("M31").Value = -2577.7058

$ws.Range("H34").Value = 4834.846
$ws.Range("I34").Value = 2872.7058
$ws.Range("K34").Value = 2872.7058
$ws.Range("M34").Value = -2670.7058

$ws.Range("H86").Value = 7429.5
$ws.Range("I86").Value = 7984
$ws.Range("J86").Value = 6875
$ws.Range("K86").Value = 7984
$ws.Range("L86").Value = 6875
$ws.Range("M86").Value = -6861
$ws.Range("N86").Value = -9121

$ws.Range("H89").Value = 7429.5
$ws.Range("I89").Value = 7984
$ws.Range("J89").Value = 6875
$ws.Range("K89").Value = 39920
$ws.Range("L89").Value = 34375
$ws.Range("M89").Value = -34304
$ws.Range("N89").Value = -45607

$ws.Range("H99").Value = 1558.5
$ws.Range("I99").Value = 1642.5714
$ws.Range("J99").Value = 1362.3334
$ws.Range("K99").Value = 1642.5714
$ws.Range("L99").Value = 1362.3334
$ws.Range("M99").Value = -144.5714
$ws.Range("N99").Value = -4358.3334

$ws.Range("H126").Value = 1558.5
$ws.Range("I126").Value = 1642.5714
$ws.Range("J126").Value = 1362.3334
$ws.Range("K126").Value = 4927.7142
$ws.Range("L126").Value = 4087.0002
$ws.Range("M126").Value = -2457.7142
$ws.Range("N126").Value = -9027.0002

$ws.Range("H134").Value = 6868
$ws.Range("I134").Value = 8097.5
$ws.Range("K134").Value = 24292.5
$ws.Range("M134").Value = -21757.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2584.818
$ws.Range("I107").Value = 4664.4
$ws.Range("K107").Value = 13993.2
$ws.Range("M107").Value = -12073.2

$ws.Range("H113").Value = 726.1429
$ws.Range("I113").Value = 632.6667
$ws.Range("J113").Value = 796.25
$ws.Range("K113").Value = 1898.0001
$ws.Range("L113").Value = 2388.75
$ws.Range("M113").Value = 271.9999
$ws.Range("N113").Value = -6728.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1394.5
$ws.Range("I113").Value = 1394.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1394.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 775.5
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 4761.154
$ws.Range("I132").Value = 4790.643
$ws.Range("K132").Value = 14371.929
$ws.Range("M132").Value = -11841.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3156.5715
$ws.Range("I40").Value = 2865.1667
$ws.Range("J40").Value = 4905
$ws.Range("K40").Value = 2865.1667
$ws.Range("L40").Value = 4905
$ws.Range("M40").Value = -2729.1667
$ws.Range("N40").Value = -5177

$ws.Range("H43").Value = 18661.666
$ws.Range("J43").Value = 18661.666
$ws.Range("L43").Value = 18661.666
$ws.Range("N43").Value = -19047.666

$ws.Range("H46").Value = 3310
$ws.Range("I46").Value = 2275
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 2275
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -2087
$ws.Range("N46").Value = -4376

$ws.Range("H82").Value = 943.7143
$ws.Range("I82").Value = 943.7143
$ws.Range("K82").Value = 943.7143
$ws.Range("M82").Value = -582.7143

$ws.Range("H85").Value = 943.7143
$ws.Range("I85").Value = 943.7143
$ws.Range("K85").Value = 943.7143
$ws.Range("M85").Value = 304.2857

$ws.Range("H132").Value = 5478.8887
$ws.Range("I132").Value = 4445.6665
$ws.Range("K132").Value = 13336.9995
$ws.Range("M132").Value = -10806.9995

$ws.Range("H136").Value = 22569.959
$ws.Range("I136").Value = 2699
$ws.Range("K136").Value = 8097
$ws.Range("M136").Value = -5547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38930
$ws.Range("J54").Value = 38930
$ws.Range("L54").Value = 38930
$ws.Range("N54").Value = -39970

$ws.Range("H81").Value = 3016.7
$ws.Range("I81").Value = 3016.7
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6033.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4972.4
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 3016.7
$ws.Range("I84").Value = 3016.7
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 30167
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -24863
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 499.5
$ws.Range("I122").Value = 499.5
$ws.Range("K122").Value = 1498.5
$ws.Range("M122").Value = 951.5

$ws.Range("H132").Value = 2698.606
$ws.Range("I132").Value = 2455.08
$ws.Range("K132").Value = 7365.24
$ws.Range("M132").Value = -4835.24
